# Generalize source data: pull the "source" fields out of SubjectAssertion
# into their own SourceData class/sheet, and rename Procedure's
# age_at_observation column to age_at_event.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new "SourceData" sheet right after "SubjectAssertion".
# ------------------------------------------------------------------
$subjectAssertion = $wb.Worksheets.Item("SubjectAssertion")
$sourceData = $wb.Worksheets.Add($null, $subjectAssertion)
$sourceData.Name = "SourceData"

$sourceDataHeaders = @(
    "code",
    "display",
    "value_code",
    "value_display",
    "value_number",
    "value_units",
    "value_units_display",
    "has_access_policy",
    "id"
)
for ($i = 0; $i -lt $sourceDataHeaders.Length; $i++) {
    $sourceData.Cells.Item(1, $i + 1).Value = $sourceDataHeaders[$i]
}

# ------------------------------------------------------------------
# 2. Re-shuffle the SubjectAssertion header row: move the age_at_*
#    columns up front, drop source_code/source_display (now modeled by
#    SourceData), and add value_units_display + source_data.
# ------------------------------------------------------------------
$subjectAssertionHeaders = @(
    "assertion_type",
    "age_at_assertion",
    "age_at_event",
    "age_at_resolution",
    "code",
    "display",
    "value_code",
    "value_display",
    "value_number",
    "value_units",
    "value_units_display",
    "source_data",
    "has_access_policy",
    "id"
)
for ($i = 0; $i -lt $subjectAssertionHeaders.Length; $i++) {
    $subjectAssertion.Cells.Item(1, $i + 1).Value = $subjectAssertionHeaders[$i]
}

# ------------------------------------------------------------------
# 3. Procedure sheet: age_at_observation -> age_at_event
# ------------------------------------------------------------------
$procedure = $wb.Worksheets.Item("Procedure")
$procedure.Range("C1").Value = "age_at_event"

Write-Output "Generalized source data: added SourceData sheet, updated SubjectAssertion and Procedure headers."
